$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphIndexStartingWith([string]$prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Insert a new "8 February 1785" entry right before the existing
#    "11 February 1785" entry (i.e. right after the "7 February 1785" entry).
# ---------------------------------------------------------------------------
$idx7Feb = Find-ParagraphIndexStartingWith("7 February 1785")
if ($idx7Feb -lt 0) { throw "Could not find '7 February 1785' paragraph" }

$anchor = $d.Paragraphs.Item($idx7Feb).Range
$anchor.InsertParagraphAfter()

$newParaIdx = $idx7Feb + 1
$newPara = $d.Paragraphs.Item($newParaIdx)

$xml8Feb = "<w:p $wNs><w:pPr><w:rPr><w:color w:val=`"000000`"/></w:rPr></w:pPr>" +
           "<w:r><w:rPr><w:b/><w:color w:val=`"000000`"/></w:rPr><w:t>8 February 1785</w:t></w:r>" +
           "<w:r><w:rPr><w:color w:val=`"000000`"/></w:rPr><w:t xml:space=`"preserve`">  Maximilian-Franz, Elector-Archbishop of Cologne, restores Christian Gottlob Neefe (37) to his full salary as court organist.  See 27 June 1784.</w:t></w:r>" +
           "</w:p>"

$newPara.Range.InsertXML($xml8Feb)

# ---------------------------------------------------------------------------
# 2) Insert a new "10 April 1785" entry right before the existing
#    "14 April 1785" entry (i.e. right after the "6 April 1785" entry).
# ---------------------------------------------------------------------------
$idx6Apr = Find-ParagraphIndexStartingWith("6 April 1785")
if ($idx6Apr -lt 0) { throw "Could not find '6 April 1785' paragraph" }

$anchor2 = $d.Paragraphs.Item($idx6Apr).Range
$anchor2.InsertParagraphAfter()

$newParaIdx2 = $idx6Apr + 1
$newPara2 = $d.Paragraphs.Item($newParaIdx2)

$xml10Apr = "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:eastAsia=`"Cambria`"/></w:rPr></w:pPr>" +
            "<w:r><w:rPr><w:rFonts w:eastAsia=`"Cambria`"/><w:b/></w:rPr><w:t>10 April 1785</w:t></w:r>" +
            "<w:r><w:rPr><w:rFonts w:eastAsia=`"Cambria`"/></w:rPr><w:t xml:space=`"preserve`">  While returning to Salem, North Carolina from the Friedland Church,  a short distance away, Johann Friedrich Peter (38) falls from his horse and is dragged some distance.  Unconscious for about 30 minutes, he is discovered by two passing men who take him back to Friedland.  He recovers and is able to return to Salem on 11 April.</w:t></w:r>" +
            "</w:p>"

$newPara2.Range.InsertXML($xml10Apr)

# ---------------------------------------------------------------------------
# 3) "4 January 2016" -> "4 July 2016", splitting "January 2016" into two
#    runs: "July" and " 2016" (the "4 " run is left untouched).
# ---------------------------------------------------------------------------
$idxDate = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt.TrimEnd([char]13) -eq "4 January 2016") {
        $idxDate = $i
        break
    }
}
if ($idxDate -lt 0) { throw "Could not find '4 January 2016' paragraph" }

$datePara = $d.Paragraphs.Item($idxDate)
$xmlDate = "<w:p $wNs w:rsidR=`"00641EF9`" w:rsidRPr=`"00371401`" w:rsidRDefault=`"00641EF9`" w:rsidP=`"00641EF9`">" +
           "<w:r w:rsidRPr=`"00371401`"><w:t xml:space=`"preserve`">4 </w:t></w:r>" +
           "<w:r><w:t>July</w:t></w:r>" +
           "<w:r><w:t xml:space=`"preserve`"> 2016</w:t></w:r></w:p>"

$datePara.Range.InsertXML($xmlDate)
